# Applies:
#  1. Table on slide 5 switches from the custom "Table_0" style to the
#     built-in "Medium Style 2 - Accent 1" table style.
#  2. The presentation's theme (currently the "Integral"/"Red Violet" theme)
#     is changed over to the stock "Office Theme"/"Office" colour scheme.

$p = $ppt.ActivePresentation

# -- 1. Table style -----------------------------------------------------
$slide = $p.Slides.Item(5)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{61D7C6ED-97F4-4B6F-8128-CC0297146A02}")
    }
}

# -- 2. Theme colours -----------------------------------------------------
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

$master = $p.SlideMaster
$scheme = $master.Theme.ThemeColorScheme
for ($i = 1; $i -le $scheme.Count; $i++) {
    $scheme.Item($i).RGB = $officeColors[$i - 1]
}
